# Hemos cambiado la fórmula de Ventas objetivo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column R ("uds. Objetivo semana pasada") values for the affected rows
$ws.Range("R4").Value = 3
$ws.Range("R8").Value = 2
$ws.Range("R9").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("R12").Value = 3
$ws.Range("R13").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("R15").Value = 5
$ws.Range("R18").Value = 2
$ws.Range("R21").Value = 2
$ws.Range("R22").Value = 1
$ws.Range("R24").Value = 1
$ws.Range("R31").Value = 1

# Update column T ("Tendencia Consumo") for rows whose trend changed as a result
$ws.Range("T10").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("T24").Value = 0
$ws.Range("T31").Value = 0

# Update column U ("Pedido Final") for row 31
$ws.Range("U31").Value = 0

# Row 31 becomes hidden
$ws.Rows(31).Hidden = $true

# Update the Total_Unidades summary cell
$ws.Range("C38").Value = 41
